$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4084957242012024
$ws.Range("B1").Value = 0.8098462820053101
$ws.Range("C1").Value = 4.161284446716309
$ws.Range("D1").Value = 2.051283836364746
$ws.Range("E1").Value = 1.07550323009491
